$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 397.73
$ws.Range("I15").Value = 397.73
$ws.Range("K15").Value = 1193.19
$ws.Range("M15").Value = -1024.19
# Row 17
$ws.Range("H17").Value = 2031068.1
$ws.Range("J17").Value = 2031068.1
$ws.Range("L17").Value = 6093204.300000001
$ws.Range("N17").Value = -6093540.300000001
# Row 76
$ws.Range("H76").Value = 4172604.2
$ws.Range("I76").Value = 5561763
$ws.Range("J76").Value = 5128.6
$ws.Range("K76").Value = 5561763
$ws.Range("L76").Value = 5128.6
$ws.Range("M76").Value = -5561448
$ws.Range("N76").Value = -5758.6
# Row 79
$ws.Range("H79").Value = 4172604.2
$ws.Range("I79").Value = 5561763
$ws.Range("J79").Value = 5128.6
$ws.Range("K79").Value = 5561763
$ws.Range("L79").Value = 5128.6
$ws.Range("M79").Value = -5560671
$ws.Range("N79").Value = -7312.6
# Row 112
$ws.Range("H112").Value = 1056.8966
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1098.0769
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 3294.2307
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -5510.2307
# Row 129
$ws.Range("H129").Value = 942.07245
$ws.Range("I129").Value = 689.36365
$ws.Range("J129").Value = 990
$ws.Range("K129").Value = 2068.09095
$ws.Range("L129").Value = 2970
$ws.Range("M129").Value = 2931.90905
$ws.Range("N129").Value = -12970
# Row 132
$ws.Range("H132").Value = 1076.1578
$ws.Range("I132").Value = 867.1613
$ws.Range("J132").Value = 2001.7142
$ws.Range("K132").Value = 2601.4839
$ws.Range("L132").Value = 6005.142599999999
$ws.Range("M132").Value = -71.48390000000018
$ws.Range("N132").Value = -11065.1426
# Row 135
$ws.Range("H135").Value = 1386.9434
$ws.Range("I135").Value = 1176.878
$ws.Range("J135").Value = 2104.6667
$ws.Range("K135").Value = 10591.902
$ws.Range("L135").Value = 18942.0003
$ws.Range("M135").Value = -8056.902
$ws.Range("N135").Value = -24012.0003
# Row 137
$ws.Range("H137").Value = 965.0682
$ws.Range("I137").Value = 792.9178000000001
$ws.Range("J137").Value = 1802.8667
$ws.Range("K137").Value = 2378.7534
$ws.Range("L137").Value = 5408.6001
$ws.Range("M137").Value = 171.2465999999999
$ws.Range("N137").Value = -10508.6001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1103.35
$ws.Range("I2").Value = 902
$ws.Range("J2").Value = 1573.1666
$ws.Range("K2").Value = 902
$ws.Range("L2").Value = 1573.1666
$ws.Range("M2").Value = -789
$ws.Range("N2").Value = -1799.1666
# Row 32
$ws.Range("H32").Value = 956.04
$ws.Range("I32").Value = 858.1064
$ws.Range("J32").Value = 2490.3333
$ws.Range("K32").Value = 858.1064
$ws.Range("L32").Value = 2490.3333
$ws.Range("M32").Value = -571.1064
$ws.Range("N32").Value = -3064.3333
# Row 61
$ws.Range("H61").Value = 2777.724
$ws.Range("I61").Value = 2944.28
$ws.Range("J61").Value = 1736.75
$ws.Range("K61").Value = 2944.28
$ws.Range("L61").Value = 1736.75
$ws.Range("M61").Value = -2732.28
$ws.Range("N61").Value = -2160.75
# Row 116
$ws.Range("H116").Value = 1103.35
$ws.Range("I116").Value = 902
$ws.Range("J116").Value = 1573.1666
$ws.Range("K116").Value = 902
$ws.Range("L116").Value = 1573.1666
$ws.Range("M116").Value = 1392
$ws.Range("N116").Value = -6161.1666
# Row 132
$ws.Range("H132").Value = 2224861.8
$ws.Range("I132").Value = 1957.7941
$ws.Range("K132").Value = 5873.3823
$ws.Range("M132").Value = -3343.3823
# Row 136
$ws.Range("H136").Value = 2777.724
$ws.Range("I136").Value = 2944.28
$ws.Range("J136").Value = 1736.75
$ws.Range("K136").Value = 8832.84
$ws.Range("L136").Value = 5210.25
$ws.Range("M136").Value = -6282.84
$ws.Range("N136").Value = -10310.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1103.35
$ws.Range("I3").Value = 902
$ws.Range("J3").Value = 1573.1666
$ws.Range("K3").Value = 902
$ws.Range("L3").Value = 1573.1666
$ws.Range("M3").Value = -788
$ws.Range("N3").Value = -1801.1666
# Row 134
$ws.Range("H134").Value = 3949.5476
$ws.Range("I134").Value = 4673
$ws.Range("J134").Value = 2335.6924
$ws.Range("K134").Value = 14019
$ws.Range("L134").Value = 7007.0772
$ws.Range("M134").Value = -11484
$ws.Range("N134").Value = -12077.0772

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5551.342
$ws.Range("I31").Value = 1476.6415
$ws.Range("K31").Value = 1476.6415
$ws.Range("M31").Value = -1181.6415
# Row 34
$ws.Range("H34").Value = 5551.342
$ws.Range("I34").Value = 1476.6415
$ws.Range("K34").Value = 1476.6415
$ws.Range("M34").Value = -1274.6415
# Row 58
$ws.Range("H58").Value = 900.8933
$ws.Range("I58").Value = 546.0566
$ws.Range("K58").Value = 546.0566
$ws.Range("M58").Value = -343.0566
# Row 122
$ws.Range("H122").Value = 2811.1333
$ws.Range("I122").Value = 2416.75
$ws.Range("J122").Value = 3261.8572
$ws.Range("K122").Value = 7250.25
$ws.Range("L122").Value = 9785.571599999999
$ws.Range("M122").Value = -4800.25
$ws.Range("N122").Value = -14685.5716
# Row 134
$ws.Range("H134").Value = 1599.8518
$ws.Range("I134").Value = 1940.3405
$ws.Range("J134").Value = 1129.1765
$ws.Range("K134").Value = 5821.0215
$ws.Range("L134").Value = 3387.5295
$ws.Range("M134").Value = -3286.0215
$ws.Range("N134").Value = -8457.529500000001
# Row 136
$ws.Range("H136").Value = 900.8933
$ws.Range("I136").Value = 546.0566
$ws.Range("K136").Value = 1638.1698
$ws.Range("M136").Value = 911.8301999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 9819731
$ws.Range("I137").Value = 7742.353
$ws.Range("J137").Value = 19631720
$ws.Range("K137").Value = 23227.059
$ws.Range("L137").Value = 58895160
$ws.Range("M137").Value = -18127.059
$ws.Range("N137").Value = -58905360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1531.8485
$ws.Range("I132").Value = 1138.9038
$ws.Range("J132").Value = 2991.3572
$ws.Range("K132").Value = 3416.7114
$ws.Range("L132").Value = 8974.071599999999
$ws.Range("M132").Value = -886.7114000000001
$ws.Range("N132").Value = -14034.0716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 60
$ws.Range("H60").Value = 33000
$ws.Range("J60").Value = 33000
$ws.Range("L60").Value = 33000
$ws.Range("N60").Value = -34018
# Row 132
$ws.Range("H132").Value = 9718407
$ws.Range("I132").Value = 13362067
$ws.Range("J132").Value = 1981.2667
$ws.Range("K132").Value = 40086201
$ws.Range("L132").Value = 5943.800099999999
$ws.Range("M132").Value = -40083671
$ws.Range("N132").Value = -11003.8001
# Row 136
$ws.Range("H136").Value = 4748.6924
$ws.Range("I136").Value = 2915.2153
$ws.Range("J136").Value = 13916.077
$ws.Range("K136").Value = 8745.6459
$ws.Range("L136").Value = 41748.231
$ws.Range("M136").Value = -6195.6459
$ws.Range("N136").Value = -46848.231

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 710.36505
$ws.Range("I132").Value = 465.1875
$ws.Range("J132").Value = 1494.9333
$ws.Range("K132").Value = 1395.5625
$ws.Range("L132").Value = 4484.7999
$ws.Range("M132").Value = 1134.4375
$ws.Range("N132").Value = -9544.7999
# Row 136
$ws.Range("H136").Value = 7464865
$ws.Range("I136").Value = 2254.875
$ws.Range("J136").Value = 26317774
$ws.Range("K136").Value = 6764.625
$ws.Range("L136").Value = 78953322
$ws.Range("M136").Value = -4214.625
$ws.Range("N136").Value = -78958422
